$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 100: new Abies koreana / 200-250 CM RB line (copy of existing product,
#     new quantity/price), formatted like the row above it (row 99) ---
$ws.Range("A99:H99").Copy()
$ws.Range("A100:H100").PasteSpecial(-4122)

$ws.Range("A100").Value = 5000
$ws.Range("B100").Value = "ABKOREAN"
$ws.Range("C100").Value = "2200250M"
$ws.Range("D100").Value = "Abies koreana"
$ws.Range("E100").Value = "200-250 CM RB"
$ws.Range("F100").Value = 100
$ws.Range("G100").Value = 43511
$ws.Range("H100").Value = 43525

# --- Row 101: brand-new product line added during the import merge,
#     no date columns filled in yet ---
$ws.Range("A99:F99").Copy()
$ws.Range("A101:F101").PasteSpecial(-4122)

$ws.Range("A101").Value = 3000
$ws.Range("C101").Value = "2200250M"
$ws.Range("D101").Value = "Barry Hebbron"
$ws.Range("E101").Value = "152CM"
$ws.Range("B101").Value = "BDHNEW"
$ws.Range("F101").Value = 2000

$ws.Range("F101").Select()
